$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Row height tweaks on the two "{{s[3]...}}" rows (both were 404,
#    first semester row -> 246, second semester row -> 283).
# ------------------------------------------------------------------
$table1 = $d.Tables.Item(1)
$table1.Rows.Item(6).Height = 246 / 20.0

$table2 = $d.Tables.Item(2)
$table2.Rows.Item(6).Height = 283 / 20.0

# ------------------------------------------------------------------
# 2) Trim the trailing school year from the header text. A throwaway
#    bookmark is dropped at the boundary with the preceding "{{name}}"
#    run first so the two runs (which already differ only by rsid,
#    not by visible formatting) are not silently coalesced back into
#    a single run by the edit.
# ------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("                           Grade 8 School Year 2025/2026", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$guard = $d.Range($anchor.Start, $anchor.Start)
$d.Bookmarks.Add("zzGuard", $guard) | Out-Null

$d.Content.Find.Execute("Grade 8 School Year 2025/2026", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Grade 8 School Year ", 2) | Out-Null

if ($d.Bookmarks.Exists("zzGuard")) {
    $d.Bookmarks.Item("zzGuard").Delete()
}

# ------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the very start of the document
#    to right after the opening "{{" of the {{s[3]["4"]}} merge tag
#    (this also removes it from its old location, since a document can
#    only have one bookmark with a given name).
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute('{{s[3]["4"]}}', $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$splitPoint = $d.Range($target.Start + 2, $target.Start + 2)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null
